# edit.ps1 - applies the changes described by the target diff:
#
#  1. The table on slide 16 (the "Total Outflow" / "Total Inflow" summary
#     table) is re-styled: its <a:tableStyleId> changes from
#     {A389863D-27FF-4290-B140-763D5B3AF9F0} to
#     {365F555D-2F36-4AA3-8DFC-432B71DFD07A}.
#
#  2. The deck's theme ("theme1.xml", the theme actually used by the
#     slide master / all slides) is swapped from the old "Integral"
#     colour palette to the "Office Theme" default colour palette
#     (which, before this edit, only lived - unused by any slide - in
#     theme2.xml, the Notes Master's theme). The two themes already
#     share an identical font scheme and format scheme, so only the
#     12 theme colours actually need to change.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style the table on slide 16.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $tableShape = $shape
    }
}
$tableShape.Table.ApplyStyle("{365F555D-2F36-4AA3-8DFC-432B71DFD07A}")

# ---------------------------------------------------------------------
# 2) Swap the theme colour scheme used by the slides (theme1.xml) from
#    "Integral" to the stock "Office Theme" palette.
#    ThemeColorScheme index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
#    5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#    11 hlink, 12 folHlink.
#    .RGB uses the standard VBA RGB() packing: R + G*256 + B*65536.
# ---------------------------------------------------------------------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

$colorScheme.Item(1).RGB  = 0 + (0 * 256) + (0 * 65536)          # dk1      000000
$colorScheme.Item(2).RGB  = 255 + (255 * 256) + (255 * 65536)    # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 68 + (84 * 256) + (106 * 65536)      # dk2      44546A
$colorScheme.Item(4).RGB  = 231 + (230 * 256) + (230 * 65536)    # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 91 + (155 * 256) + (213 * 65536)     # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 237 + (125 * 256) + (49 * 65536)     # accent2  ED7D31
$colorScheme.Item(7).RGB  = 165 + (165 * 256) + (165 * 65536)    # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 255 + (192 * 256) + (0 * 65536)      # accent4  FFC000
$colorScheme.Item(9).RGB  = 68 + (114 * 256) + (196 * 65536)     # accent5  4472C4
$colorScheme.Item(10).RGB = 112 + (173 * 256) + (71 * 65536)     # accent6  70AD47
$colorScheme.Item(11).RGB = 5 + (99 * 256) + (193 * 65536)       # hlink    0563C1
$colorScheme.Item(12).RGB = 149 + (79 * 256) + (114 * 65536)     # folHlink 954F72
